$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14263
$ws1.Range("F3").Value = 335
$ws1.Range("F4").Value = 689
$ws1.Range("F6").Value = 556
$ws1.Range("F7").Value = 1491

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14263
$ws4.Range("F3").Value = 335
$ws4.Range("F4").Value = 689
$ws4.Range("F8").Value = 556
$ws4.Range("F9").Value = 1491
